$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 140; all rows below shift up by one, matching the
# correction described in the commit message ("Correcting errors in
# trial matrix and code").
$ws.Rows.Item(140).Delete()

# Update the view so the active cell / top-left cell match the target state.
$ws.Application.ActiveWindow.ScrollRow = 121
$ws.Range("A140").Select()
